$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F9").Value = 1
$ws.Range("F11").Value = -1
$ws.Range("F14").Value = -1
$ws.Range("F16").Value = 2
$ws.Range("F19").Value = -2
$ws.Range("F23").Value = 0
$ws.Range("F25").Value = -4
